# Applies updated market-price / profit figures to the Leve profit sheets.
# Values below come from a scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1577.8959
$ws.Range("I17").Value = 1263
$ws.Range("J17").Value = 1598.8889
$ws.Range("K17").Value = 3789
$ws.Range("L17").Value = 4796.6667
$ws.Range("M17").Value = -3621
$ws.Range("N17").Value = -5132.6667
$ws.Range("H19").Value = 55556610
$ws.Range("J19").Value = 62501036
$ws.Range("L19").Value = 62501036
$ws.Range("N19").Value = -62501386
$ws.Range("H43").Value = 429
$ws.Range("I43").Value = 390
$ws.Range("J43").Value = 456.85715
$ws.Range("K43").Value = 390
$ws.Range("L43").Value = 456.85715
$ws.Range("M43").Value = -321
$ws.Range("N43").Value = -594.85715
$ws.Range("H51").Value = 10105909
$ws.Range("I51").Value = 22731648
$ws.Range("J51").Value = 5318
$ws.Range("K51").Value = 22731648
$ws.Range("L51").Value = 5318
$ws.Range("M51").Value = -22731164
$ws.Range("N51").Value = -6286
$ws.Range("H53").Value = 392.1
$ws.Range("I53").Value = 163.07143
$ws.Range("J53").Value = 926.5
$ws.Range("K53").Value = 163.07143
$ws.Range("L53").Value = 926.5
$ws.Range("M53").Value = 473.92857
$ws.Range("N53").Value = -2200.5
$ws.Range("H55").Value = 363.5263
$ws.Range("I55").Value = 366.5
$ws.Range("J55").Value = 360.22223
$ws.Range("K55").Value = 366.5
$ws.Range("L55").Value = 360.22223
$ws.Range("M55").Value = -152.5
$ws.Range("N55").Value = -788.2222300000001
$ws.Range("H76").Value = 2870.55
$ws.Range("I76").Value = 2601.125
$ws.Range("J76").Value = 3948.25
$ws.Range("K76").Value = 2601.125
$ws.Range("L76").Value = 3948.25
$ws.Range("M76").Value = -2286.125
$ws.Range("N76").Value = -4578.25
$ws.Range("H79").Value = 2870.55
$ws.Range("I79").Value = 2601.125
$ws.Range("J79").Value = 3948.25
$ws.Range("K79").Value = 2601.125
$ws.Range("L79").Value = 3948.25
$ws.Range("M79").Value = -1509.125
$ws.Range("N79").Value = -6132.25
$ws.Range("H94").Value = 2401
$ws.Range("I94").Value = 2401
$ws.Range("K94").Value = 2401
$ws.Range("M94").Value = -1950
$ws.Range("H98").Value = 1920.875
$ws.Range("I98").Value = 2022.2778
$ws.Range("J98").Value = 1616.6666
$ws.Range("K98").Value = 2022.2778
$ws.Range("L98").Value = 1616.6666
$ws.Range("M98").Value = -524.2778000000001
$ws.Range("N98").Value = -4612.6666
$ws.Range("H122").Value = 1920.875
$ws.Range("I122").Value = 2022.2778
$ws.Range("J122").Value = 1616.6666
$ws.Range("K122").Value = 6066.8334
$ws.Range("L122").Value = 4849.9998
$ws.Range("M122").Value = -3616.8334
$ws.Range("N122").Value = -9749.9998
$ws.Range("H129").Value = 1163.7361
$ws.Range("I129").Value = 879.7778
$ws.Range("J129").Value = 1204.3016
$ws.Range("K129").Value = 2639.3334
$ws.Range("L129").Value = 3612.9048
$ws.Range("M129").Value = 2360.6666
$ws.Range("N129").Value = -13612.9048
$ws.Range("H137").Value = 1820.2188
$ws.Range("I137").Value = 3096.4167
$ws.Range("K137").Value = 9289.250100000001
$ws.Range("M137").Value = -6739.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 4000
$ws.Range("K57").Value = 4000
$ws.Range("M57").Value = -3516
$ws.Range("H63").Value = 3687.6667
$ws.Range("I63").Value = 2151.25
$ws.Range("J63").Value = 9833.333000000001
$ws.Range("K63").Value = 2151.25
$ws.Range("L63").Value = 9833.333000000001
$ws.Range("M63").Value = -1465.25
$ws.Range("N63").Value = -11205.333
$ws.Range("H66").Value = 3687.6667
$ws.Range("I66").Value = 2151.25
$ws.Range("J66").Value = 9833.333000000001
$ws.Range("K66").Value = 10756.25
$ws.Range("L66").Value = 49166.665
$ws.Range("M66").Value = -7324.25
$ws.Range("N66").Value = -56030.665
$ws.Range("H97").Value = 6747383.5
$ws.Range("I97").Value = 1042934.75
$ws.Range("J97").Value = 25001620
$ws.Range("K97").Value = 1042934.75
$ws.Range("L97").Value = 25001620
$ws.Range("M97").Value = -1042438.75
$ws.Range("N97").Value = -25002612

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1893.65
$ws.Range("I105").Value = 1507.4615
$ws.Range("J105").Value = 2610.8572
$ws.Range("K105").Value = 1507.4615
$ws.Range("L105").Value = 2610.8572
$ws.Range("M105").Value = 239.5385000000001
$ws.Range("N105").Value = -6104.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4275493.5
$ws.Range("I31").Value = 1193.9783
$ws.Range("J31").Value = 10419799
$ws.Range("K31").Value = 1193.9783
$ws.Range("L31").Value = 10419799
$ws.Range("M31").Value = -898.9783
$ws.Range("N31").Value = -10420389
$ws.Range("H34").Value = 4275493.5
$ws.Range("I34").Value = 1193.9783
$ws.Range("J34").Value = 10419799
$ws.Range("K34").Value = 1193.9783
$ws.Range("L34").Value = 10419799
$ws.Range("M34").Value = -991.9783
$ws.Range("N34").Value = -10420203
$ws.Range("H53").Value = 31833.334
$ws.Range("J53").Value = 31833.334
$ws.Range("L53").Value = 31833.334
$ws.Range("N53").Value = -33047.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 337.7857
$ws.Range("J33").Value = 394
$ws.Range("L33").Value = 2364
$ws.Range("N33").Value = -2930
$ws.Range("H105").Value = 9543.637000000001
$ws.Range("J105").Value = 9995.4
$ws.Range("L105").Value = 29986.2
$ws.Range("N105").Value = -35228.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 90006
$ws.Range("J20").Value = 90006
$ws.Range("L20").Value = 90006
$ws.Range("N20").Value = -90496
$ws.Range("H70").Value = 5275.75
$ws.Range("I70").Value = 4799.1665
$ws.Range("J70").Value = 5633.1875
$ws.Range("K70").Value = 4799.1665
$ws.Range("L70").Value = 5633.1875
$ws.Range("M70").Value = -4529.1665
$ws.Range("N70").Value = -6173.1875
$ws.Range("H73").Value = 5275.75
$ws.Range("I73").Value = 4799.1665
$ws.Range("J73").Value = 5633.1875
$ws.Range("K73").Value = 4799.1665
$ws.Range("L73").Value = 5633.1875
$ws.Range("M73").Value = -3863.1665
$ws.Range("N73").Value = -7505.1875
$ws.Range("H80").Value = 2580.111
$ws.Range("I80").Value = 2663
$ws.Range("J80").Value = 2476.5
$ws.Range("K80").Value = 2663
$ws.Range("L80").Value = 2476.5
$ws.Range("M80").Value = -1665
$ws.Range("N80").Value = -4472.5
$ws.Range("H83").Value = 2580.111
$ws.Range("I83").Value = 2663
$ws.Range("J83").Value = 2476.5
$ws.Range("K83").Value = 13315
$ws.Range("L83").Value = 12382.5
$ws.Range("M83").Value = -8323
$ws.Range("N83").Value = -22366.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5200.4
$ws.Range("I61").Value = 5444.8887
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 5444.8887
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -5242.8887
$ws.Range("N61").Value = -3404
$ws.Range("H68").Value = 1531.6976
$ws.Range("I68").Value = 1448.75
$ws.Range("J68").Value = 1773
$ws.Range("K68").Value = 1448.75
$ws.Range("L68").Value = 1773
$ws.Range("M68").Value = -699.75
$ws.Range("N68").Value = -3271
$ws.Range("H71").Value = 1531.6976
$ws.Range("I71").Value = 1448.75
$ws.Range("J71").Value = 1773
$ws.Range("K71").Value = 7243.75
$ws.Range("L71").Value = 8865
$ws.Range("M71").Value = -3499.75
$ws.Range("N71").Value = -16353
$ws.Range("H113").Value = 5200.4
$ws.Range("I113").Value = 5444.8887
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 5444.8887
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -3274.8887
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 318
$ws.Range("I14").Value = 202
$ws.Range("J14").Value = 550
$ws.Range("K14").Value = 202
$ws.Range("L14").Value = 550
$ws.Range("M14").Value = -34
$ws.Range("N14").Value = -886
$ws.Range("H81").Value = 2306.75
$ws.Range("I81").Value = 2600.1428
$ws.Range("J81").Value = 1896
$ws.Range("K81").Value = 5200.2856
$ws.Range("L81").Value = 3792
$ws.Range("M81").Value = -4139.2856
$ws.Range("N81").Value = -5914
$ws.Range("H84").Value = 2306.75
$ws.Range("I84").Value = 2600.1428
$ws.Range("J84").Value = 1896
$ws.Range("K84").Value = 26001.428
$ws.Range("L84").Value = 18960
$ws.Range("M84").Value = -20697.428
$ws.Range("N84").Value = -29568
